$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data-entry error being fixed: the "source" (K2), "link" (L2) and "notes" (M2)
# cells for the earnings_effect_high_vs_low row had referenced the wrong paper link.
# The link cell also carried a live hyperlink to the old (wrong) URL, which is no
# longer wanted now that the link text itself is being corrected.
$ws.Range("L2").Hyperlinks.Delete()

$ws.Range("K2").Value = "Pawlowski et al. (2019) Figure 2 & Text p. 19"
$ws.Range("L2").Value = "https://ideas.repec.org/p/usg/econwp/201906.html"
$ws.Range("M2").Value = "The authors do not present their estimates in a table with standard errors. Instead they show the effect for each year in Figure 2. The effect is also somewhat jumpy between years.  However, the authors mention in the text what they believe to be the average effect. (i.e. 260€ for men / 0 for women when comparing high to low expenditure). Looking at the graph these appear to be significant at at about 10%."

# Reflect the cursor position recorded in the saved workbook.
$ws.Range("L14").Select()
